$d = $word.ActiveDocument
$q = [char]34

# --------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph. It has a bold "Meta description" run
#    followed by a normal run with the rest of the sentence.
# --------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = $d.Styles.Item("Normal")

$label = "Meta description"
$rest = ": Experience the exciting gameplay features of 9 Blazing Cashpots Megaways. Play for free with this review, featuring a Megaways mechanic and Cashpot Feature."
$metaPara.Range.Text = $label + $rest

$boldStart = $metaPara.Range.Start
$boldEnd = $boldStart + $label.Length
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Bold = 1

# --------------------------------------------------------------------
# 2. Remove the old bold "Play 9 Blazing Cashpots Megaways..." heading
#    paragraph that used to sit just before the final (italic) paragraph.
# --------------------------------------------------------------------
$count = $d.Paragraphs.Count
$oldHeadingPara = $d.Paragraphs.Item($count - 1)
$oldHeadingPara.Range.Delete()

# --------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph with the new
#    image-generation prompt text, keeping its italic run formatting.
# --------------------------------------------------------------------
$lastParaObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$s = $lastParaObj.Range.Start
$e = $lastParaObj.Range.End
$finalRange = $d.Range($s, $e)

$newText = "Create a feature image fitting the game " + $q + "9 Blazing Cashpots Megaways" + $q + " with the following requirements: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses The feature image should be a fun and exciting depiction of the Maya warrior enjoying the game of 9 Blazing Cashpots Megaways. The warrior should be shown with a big smile on their face, wearing cool sunglasses to show off their winning streak. They could be sitting in front of the reels with a stack of colorful fruit symbols on one side and cashpot symbols on the other. The background of the image could be a mix of jungle foliage and casino elements. The overall style of the image should be cartoonish, bright, and energetic, reflecting the upbeat nature of the game."

$finalRange.Text = $newText
